$d = $word.ActiveDocument

$replacements = @(
    @("66×42=2772", "66×89=5874"),
    @("87×97=8439", "77×72=5544"),
    @("84×34=2856", "48×73=3504"),
    @("26×43=1118", "67×80=5360"),
    @("59×33=1947", "31×47=1457"),
    @("27×96=2592", "41×24=984"),
    @("95×31=2945", "51×65=3315"),
    @("90×11=990",  "60×64=3840"),
    @("88×13=1144", "64×81=5184"),
    @("68×41=2788", "50×84=4200"),
    @("81×62=5022", "94×37=3478"),
    @("80×84=6720", "19×65=1235"),
    @("82×99=8118", "36×61=2196"),
    @("94×55=5170", "56×28=1568"),
    @("75×64=4800", "74×96=7104"),
    @("16×78=1248", "21×96=2016"),
    @("85×12=1020", "26×63=1638"),
    @("64×97=6208", "69×34=2346"),
    @("93×47=4371", "56×67=3752"),
    @("89×44=3916", "46×14=644"),
    @("74×94=6956", "11×54=594"),
    @("32×90=2880", "69×15=1035"),
    @("33×88=2904", "81×33=2673"),
    @("55×26=1430", "41×13=533"),
    @("36×82=2952", "52×55=2860")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
